$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108
$yellow = 65535

# --- Row 1 headers: new HomeWork columns H03 (D) and H05 (E) ---
$ws.Range("D1").Value = "H03"
$ws.Range("E1").Value = "H05"

# --- Rows 2-6 (style 13 rows): fill in D/E scores ---
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 9.5

$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 10

$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 8.5

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 9

$ws.Range("D6").Value = 9.75
$ws.Range("E6").Value = 10

# --- Row 7: C7 also changes (0 -> 9), plus D7/E7 ---
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 7.5
$ws.Range("E7").Value = 9

# --- Rows 8-10 (style 4 rows, cells already exist blank) ---
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0

$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 9.75

$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 8

# --- Rows 11-16: D/E cells don't exist yet, need to create with style matching column C (s=4) ---
$ws.Range("C11").Copy()
$ws.Range("D11:E11").PasteSpecial($xlPasteFormats)
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 8

$ws.Range("C12").Copy()
$ws.Range("D12:E12").PasteSpecial($xlPasteFormats)
$ws.Range("D12").Value = 9.5
$ws.Range("E12").Value = 9.5

# Row 13: D13 gets a formula (half credit) highlighted in yellow, E13 normal
$ws.Range("C13").Copy()
$ws.Range("D13:E13").PasteSpecial($xlPasteFormats)
$ws.Range("D13").Formula = "=9.5/2"
$ws.Range("D13").Interior.Color = $yellow
$ws.Range("D13").HorizontalAlignment = $xlCenter
$ws.Range("E13").Value = 10

$ws.Range("C14").Copy()
$ws.Range("D14:E14").PasteSpecial($xlPasteFormats)
$ws.Range("D14").Value = 7
$ws.Range("E14").Value = 6.5

# Row 15: D15 gets a formula (half credit) highlighted in yellow (General format), E15 normal
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("D15").Interior.Color = $yellow
$ws.Range("D15").HorizontalAlignment = $xlCenter
$ws.Range("D15").Formula = "=10/2"
$ws.Range("E15").Value = 10

$ws.Range("C16").Copy()
$ws.Range("D16:E16").PasteSpecial($xlPasteFormats)
$ws.Range("D16").Value = 8.5
$ws.Range("E16").Value = 9

# --- Extra point manual edits in the summary table (C column 0 -> 1) ---
$ws.Range("C23").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C31").Value = 1

# --- Update selection to match final cursor position ---
$ws.Range("E11").Select()
